$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "MISS: $find"
    }
}

# 1. Move the period inside the closing quotation mark after "D12138"
Replace-Text "D12138”. " "D12138.” "

# 2. "future? Have to admit that Artificial" -> "future? We have to admit that artificial"
Replace-Text "future? Have to admit that Artificial" "future? We have to admit that artificial"

# 3. Remove the straight quotes around "people" in "being ""people"" through learning"
Replace-Text "being ""people"" through" "being people through"

# 4. ", etc.？" -> "? " (drop ", etc." and turn full-width ？ into ASCII "? ")
Replace-Text ", etc.？" "? "

# 5. "more expensive" -> "higher" in the high-cost investment sentence
Replace-Text "cost it takes is more expensive than traditional" "cost it takes is higher than traditional"

# 6. drop the comma after "AI-based education" before "and is not for every child"
Replace-Text "increase on AI-based education, and is not for every child" "increase on AI-based education and is not for every child"

# 7. "human interaction that requires real human skills" -> "mankind interaction that requires human skills"
Replace-Text "complex human interaction that requires real human skills" "complex mankind interaction that requires human skills"

# 8. Move the (hidden) _GoBack bookmark from the trailing empty paragraph to sit
#    between "AI-based " and "education" in "...increase on AI-based education and..."
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()
$anchor = $d.Content
$anchor.Find.Execute("increase on AI-based ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($anchor.End, $anchor.End))

Write-Output "done"
